# Adds the 2023-12-17 daily update to the violent-crime-full-year workbook.
# Each sheet (Citywide Totals, By Neighborhood, and every per-neighborhood
# sheet) keeps a running year-to-date total for 2023 in column J (plus a
# couple of small 2022 / column-I corrections). This bumps the affected
# cells to the new cumulative counts after the new day of data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 7367
$ws.Range("J3").Value = 7756
$ws.Range("I4").Value = 1776
$ws.Range("J4").Value = 1691
$ws.Range("J5").Value = 607
$ws.Range("J6").Value = 10596
$ws.Range("I7").Value = 26231
$ws.Range("J7").Value = 28017

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J2").Value = 223
$ws.Range("J7").Value = 800
$ws.Range("J8").Value = 1767
$ws.Range("J14").Value = 148
$ws.Range("J15").Value = 346
$ws.Range("J19").Value = 807
$ws.Range("J20").Value = 602
$ws.Range("J21").Value = 80
$ws.Range("J23").Value = 257
$ws.Range("J24").Value = 94
$ws.Range("J29").Value = 1492
$ws.Range("J31").Value = 295
$ws.Range("J33").Value = 1271
$ws.Range("J34").Value = 130
$ws.Range("J35").Value = 33
$ws.Range("J36").Value = 380
$ws.Range("J37").Value = 859
$ws.Range("J41").Value = 210
$ws.Range("J42").Value = 1193
$ws.Range("J43").Value = 236
$ws.Range("J46").Value = 93
$ws.Range("J48").Value = 311
$ws.Range("J50").Value = 168
$ws.Range("J52").Value = 713
$ws.Range("J54").Value = 554
$ws.Range("J57").Value = 133
$ws.Range("I63").Value = 185
$ws.Range("J63").Value = 90
$ws.Range("J65").Value = 705
$ws.Range("J67").Value = 1028
$ws.Range("J73").Value = 272
$ws.Range("J76").Value = 398
$ws.Range("J78").Value = 325
$ws.Range("J79").Value = 767
$ws.Range("J80").Value = 52
$ws.Range("J83").Value = 562
$ws.Range("J85").Value = 1150
$ws.Range("J87").Value = 94
$ws.Range("J89").Value = 348
$ws.Range("J94").Value = 313
$ws.Range("J95").Value = 401
$ws.Range("J96").Value = 311
$ws.Range("J98").Value = 207
$ws.Range("I101").Value = 26231
$ws.Range("J101").Value = 28017

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("J4").Value = 8
$ws.Range("J7").Value = 148

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J6").Value = 117
$ws.Range("J7").Value = 311

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J6").Value = 256
$ws.Range("J7").Value = 800

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J2").Value = 105
$ws.Range("J4").Value = 35
$ws.Range("J7").Value = 348

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 306
$ws.Range("J3").Value = 415
$ws.Range("J6").Value = 328
$ws.Range("J7").Value = 1150

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J3").Value = 199
$ws.Range("J7").Value = 713

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 465
$ws.Range("J3").Value = 510
$ws.Range("J6").Value = 655
$ws.Range("J7").Value = 1767

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J2").Value = 166
$ws.Range("J3").Value = 208
$ws.Range("J7").Value = 562

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 288
$ws.Range("J3").Value = 423
$ws.Range("J6").Value = 452
$ws.Range("J7").Value = 1271

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J2").Value = 141
$ws.Range("J7").Value = 401

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J5").Value = 31
$ws.Range("J6").Value = 251
$ws.Range("J7").Value = 859

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J2").Value = 203
$ws.Range("J6").Value = 265
$ws.Range("J7").Value = 705

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J2").Value = 96
$ws.Range("J3").Value = 72
$ws.Range("J6").Value = 106
$ws.Range("J7").Value = 295

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J3").Value = 388
$ws.Range("J7").Value = 1028

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J3").Value = 112
$ws.Range("J6").Value = 256
$ws.Range("J7").Value = 554

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 454
$ws.Range("J7").Value = 1492

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J6").Value = 149
$ws.Range("J7").Value = 311

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J6").Value = 313
$ws.Range("J7").Value = 807

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J6").Value = 205
$ws.Range("J7").Value = 398

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("J3").Value = 32
$ws.Range("J7").Value = 210

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J2").Value = 249
$ws.Range("J7").Value = 1193

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J3").Value = 99
$ws.Range("J6").Value = 102
$ws.Range("J7").Value = 325

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("J2").Value = 32
$ws.Range("J7").Value = 94

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("J6").Value = 39
$ws.Range("J7").Value = 93

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J6").Value = 72
$ws.Range("J7").Value = 257

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("J6").Value = 54
$ws.Range("J7").Value = 80

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 218
$ws.Range("J7").Value = 767

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J2").Value = 167
$ws.Range("J7").Value = 602

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J4").Value = 16
$ws.Range("J6").Value = 112
$ws.Range("J7").Value = 380

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("J2").Value = 36
$ws.Range("J3").Value = 34
$ws.Range("J7").Value = 130

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J3").Value = 60
$ws.Range("J7").Value = 313

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J3").Value = 70
$ws.Range("J6").Value = 160
$ws.Range("J7").Value = 346

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("J4").Value = 10
$ws.Range("J7").Value = 207

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("J3").Value = 44
$ws.Range("J7").Value = 168

$ws = $wb.Worksheets.Item("Gold Coast")
$ws.Range("J6").Value = 21
$ws.Range("J7").Value = 33

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J3").Value = 69
$ws.Range("J7").Value = 272

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("J4").Value = 18
$ws.Range("J6").Value = 81
$ws.Range("J7").Value = 223

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("J3").Value = 30
$ws.Range("J6").Value = 60
$ws.Range("J7").Value = 133

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("J6").Value = 141
$ws.Range("J7").Value = 236

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("J6").Value = 28
$ws.Range("J7").Value = 52

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("J4").Value = 9
$ws.Range("J7").Value = 94
